$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 06:49:22"

$ws1.Cells.Item(3, 1).Value = "Total filas: 47"

$ws1.Cells.Item(16, 1).Value = "05:44:02"
$ws1.Cells.Item(16, 2).Value = "06:40"
$ws1.Cells.Item(16, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(16, 4).Value = 56
$ws1.Cells.Item(16, 5).Value = "LP1912"

$ws1.Cells.Item(17, 1).Value = "06:38:54"
$ws1.Cells.Item(17, 2).Value = "06:40"
$ws1.Cells.Item(17, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(17, 4).Value = 2
$ws1.Cells.Item(17, 5).Value = "LP1912"

$ws1.Cells.Item(20, 1).Value = "06:49:22"
$ws1.Cells.Item(20, 2).Value = "06:57"
$ws1.Cells.Item(20, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(20, 4).Value = 8
$ws1.Cells.Item(20, 5).Value = "LP1912"

$ws1.Cells.Item(22, 1).Value = "06:49:22"
$ws1.Cells.Item(22, 2).Value = "06:59"
$ws1.Cells.Item(22, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(22, 4).Value = 10
$ws1.Cells.Item(22, 5).Value = "LP1912"

$ws1.Cells.Item(24, 1).Value = "06:49:22"
$ws1.Cells.Item(24, 2).Value = "07:16"
$ws1.Cells.Item(24, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(24, 4).Value = 27
$ws1.Cells.Item(24, 5).Value = "LP1912"

$ws1.Cells.Item(26, 1).Value = "06:49:22"
$ws1.Cells.Item(26, 2).Value = "07:19"
$ws1.Cells.Item(26, 3).Value = "14_ABASTO"
$ws1.Cells.Item(26, 4).Value = 30
$ws1.Cells.Item(26, 5).Value = "LP1912"

$ws1.Cells.Item(28, 1).Value = "06:49:22"
$ws1.Cells.Item(28, 2).Value = "07:21"
$ws1.Cells.Item(28, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(28, 4).Value = 32
$ws1.Cells.Item(28, 5).Value = "LP1912"

$ws1.Cells.Item(29, 1).Value = "06:49:22"
$ws1.Cells.Item(29, 2).Value = "07:21"
$ws1.Cells.Item(29, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(29, 4).Value = 32
$ws1.Cells.Item(29, 5).Value = "LP1912"

$ws1.Cells.Item(31, 1).Value = "06:49:22"
$ws1.Cells.Item(31, 2).Value = "07:29"
$ws1.Cells.Item(31, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(31, 4).Value = 40
$ws1.Cells.Item(31, 5).Value = "LP1912"

$ws1.Cells.Item(33, 1).Value = "06:49:22"
$ws1.Cells.Item(33, 2).Value = "07:35"
$ws1.Cells.Item(33, 3).Value = "10_OLMOS"
$ws1.Cells.Item(33, 4).Value = 46
$ws1.Cells.Item(33, 5).Value = "LP1912"

$ws1.Cells.Item(35, 1).Value = "06:49:22"
$ws1.Cells.Item(35, 2).Value = "07:37"
$ws1.Cells.Item(35, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(35, 4).Value = 48
$ws1.Cells.Item(35, 5).Value = "LP1912"

$ws1.Cells.Item(37, 1).Value = "06:49:22"
$ws1.Cells.Item(37, 2).Value = "07:44"
$ws1.Cells.Item(37, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(37, 4).Value = 55
$ws1.Cells.Item(37, 5).Value = "LP1912"

$ws1.Cells.Item(38, 1).Value = "06:38:54"
$ws1.Cells.Item(38, 2).Value = "07:54"
$ws1.Cells.Item(38, 3).Value = "14_ABASTO"
$ws1.Cells.Item(38, 4).Value = 76
$ws1.Cells.Item(38, 5).Value = "LP1912"

$ws1.Cells.Item(39, 1).Value = "06:49:22"
$ws1.Cells.Item(39, 2).Value = "07:55"
$ws1.Cells.Item(39, 3).Value = "14_ABASTO"
$ws1.Cells.Item(39, 4).Value = 66
$ws1.Cells.Item(39, 5).Value = "LP1912"

$ws1.Cells.Item(40, 1).Value = "06:19:59"
$ws1.Cells.Item(40, 2).Value = "07:59"
$ws1.Cells.Item(40, 3).Value = "17_ROMERO"
$ws1.Cells.Item(40, 4).Value = 100
$ws1.Cells.Item(40, 5).Value = "LP1912"

$ws1.Cells.Item(42, 1).Value = "06:49:22"
$ws1.Cells.Item(42, 2).Value = "08:00"
$ws1.Cells.Item(42, 3).Value = "17_ROMERO"
$ws1.Cells.Item(42, 4).Value = 71
$ws1.Cells.Item(42, 5).Value = "LP1912"

$ws1.Cells.Item(43, 1).Value = "06:49:22"
$ws1.Cells.Item(43, 2).Value = "08:01"
$ws1.Cells.Item(43, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(43, 4).Value = 72
$ws1.Cells.Item(43, 5).Value = "LP1912"

$ws1.Cells.Item(44, 1).Value = "06:49:22"
$ws1.Cells.Item(44, 2).Value = "08:06"
$ws1.Cells.Item(44, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(44, 4).Value = 77
$ws1.Cells.Item(44, 5).Value = "LP1912"

$ws1.Cells.Item(45, 1).Value = "06:49:22"
$ws1.Cells.Item(45, 2).Value = "08:11"
$ws1.Cells.Item(45, 3).Value = "10_OLMOS"
$ws1.Cells.Item(45, 4).Value = 82
$ws1.Cells.Item(45, 5).Value = "LP1912"

$ws1.Cells.Item(46, 1).Value = "06:19:59"
$ws1.Cells.Item(46, 2).Value = "08:12"
$ws1.Cells.Item(46, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(46, 4).Value = 113
$ws1.Cells.Item(46, 5).Value = "LP1912"

$ws1.Cells.Item(47, 1).Value = "06:49:22"
$ws1.Cells.Item(47, 2).Value = "08:13"
$ws1.Cells.Item(47, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(47, 4).Value = 84
$ws1.Cells.Item(47, 5).Value = "LP1912"

$ws1.Cells.Item(48, 1).Value = "06:38:54"
$ws1.Cells.Item(48, 2).Value = "08:28"
$ws1.Cells.Item(48, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(48, 4).Value = 110
$ws1.Cells.Item(48, 5).Value = "LP1912"

$ws1.Cells.Item(49, 1).Value = "06:49:22"
$ws1.Cells.Item(49, 2).Value = "08:29"
$ws1.Cells.Item(49, 3).Value = "15_ABASTO"
$ws1.Cells.Item(49, 4).Value = 100
$ws1.Cells.Item(49, 5).Value = "LP1912"

$ws1.Cells.Item(50, 1).Value = "06:49:22"
$ws1.Cells.Item(50, 2).Value = "08:29"
$ws1.Cells.Item(50, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(50, 4).Value = 100
$ws1.Cells.Item(50, 5).Value = "LP1912"

$ws1.Cells.Item(51, 1).Value = "06:49:22"
$ws1.Cells.Item(51, 2).Value = "08:41"
$ws1.Cells.Item(51, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(51, 4).Value = 112
$ws1.Cells.Item(51, 5).Value = "LP1912"

$ws1.Cells.Item(52, 1).Value = "06:49:22"
$ws1.Cells.Item(52, 2).Value = "08:44"
$ws1.Cells.Item(52, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(52, 4).Value = 115
$ws1.Cells.Item(52, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 06:49:22"

$ws2.Cells.Item(3, 1).Value = "Total filas: 8"

$ws2.Cells.Item(8, 1).Value = "06:49:22"
$ws2.Cells.Item(8, 2).Value = "06:57"
$ws2.Cells.Item(8, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(8, 4).Value = 8
$ws2.Cells.Item(8, 5).Value = "LP1912"

$ws2.Cells.Item(10, 1).Value = "06:49:22"
$ws2.Cells.Item(10, 2).Value = "07:16"
$ws2.Cells.Item(10, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(10, 4).Value = 27
$ws2.Cells.Item(10, 5).Value = "LP1912"

$ws2.Cells.Item(12, 1).Value = "06:49:22"
$ws2.Cells.Item(12, 2).Value = "07:44"
$ws2.Cells.Item(12, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(12, 4).Value = 55
$ws2.Cells.Item(12, 5).Value = "LP1912"

$ws2.Cells.Item(13, 1).Value = "06:49:22"
$ws2.Cells.Item(13, 2).Value = "08:44"
$ws2.Cells.Item(13, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(13, 4).Value = 115
$ws2.Cells.Item(13, 5).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 06:49:22"

$ws3.Cells.Item(3, 1).Value = "Total filas: 4"

$ws3.Cells.Item(7, 1).Value = "06:49:22"
$ws3.Cells.Item(7, 2).Value = "07:43"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 54
$ws3.Cells.Item(7, 5).Value = "L6173"

$ws3.Cells.Item(9, 1).Value = "06:49:22"
$ws3.Cells.Item(9, 2).Value = "08:36"
$ws3.Cells.Item(9, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(9, 4).Value = 107
$ws3.Cells.Item(9, 5).Value = "L6173"

